$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1)
$ws.Range("A1").Value = "Covariate"
$ws.Range("B1").Value = "Median (5%, 95% quantiles)"
$ws.Range("C1").Value = "Mean (min, max)"
$ws.Range("D1").Value = "Description"
$ws.Range("E1").Value = "UglyName"

# Update data rows 2-37 with new covariate summary data
# Row 2: cfi_site
$ws.Range("A2").Value = "CFI"
$ws.Range("B2").Value = "0.05 (0.05–0.23)"
$ws.Range("C2").Value = "0.07 (0–1)"
$ws.Range("E2").Value = "cfi_site"

# Row 3: cfi_site_with_harvest
$ws.Range("A3").Value = "CFI"
$ws.Range("B3").Value = "0.06 (0.06–0.28)"
$ws.Range("C3").Value = "0.09 (0–1)"
$ws.Range("E3").Value = "cfi_site_with_harvest"

# Row 4: cfi_site_with_vegedges
$ws.Range("A4").Value = "CFI"
$ws.Range("B4").Value = "0.06 (0.06–0.26)"
$ws.Range("C4").Value = "0.08 (0–1)"
$ws.Range("E4").Value = "cfi_site_with_vegedges"

# Row 5: fire_0_15
$ws.Range("A5").Value = "Fire <15"
$ws.Range("B5").Value = "0 (0–0.74)"
$ws.Range("C5").Value = "0.07 (0–1)"
$ws.Range("E5").Value = "fire_0_15"

# Row 6: harvest_0_15
$ws.Range("A6").Value = "Harvest <15"
$ws.Range("B6").Value = "0 (0–0.12)"
$ws.Range("C6").Value = "0.02 (0–1)"
$ws.Range("E6").Value = "harvest_0_15"

# Row 7: harvest_total
$ws.Range("A7").Value = "Harvest"
$ws.Range("B7").Value = "0 (0–0.28)"
$ws.Range("C7").Value = "0.05 (0–1)"
$ws.Range("E7").Value = "harvest_total"

# Row 8: landscape_cai_mn
$ws.Range("A8").Value = "Core Area Index"
$ws.Range("B8").Value = "0 (0–0)"
$ws.Range("C8").Value = "0 (0–0)"
$ws.Range("E8").Value = "landscape_cai_mn"

# Row 9: landscape_contag
$ws.Range("A9").Value = "Contagion"
$ws.Range("B9").Value = "62.84 (62.84–80.38)"
$ws.Range("C9").Value = "63.54 (18.38–99.53)"
$ws.Range("E9").Value = "landscape_contag"

# Row 10: landscape_ed
$ws.Range("A10").Value = "Edge Density"
$ws.Range("B10").Value = "0 (0–0)"
$ws.Range("C10").Value = "0 (0–0)"
$ws.Range("E10").Value = "landscape_ed"

# Row 11: landscape_mesh
$ws.Range("A11").Value = "Mesh Index"
$ws.Range("B11").Value = "113.08 (113.08–651.94)"
$ws.Range("C11").Value = "192.78 (0.15–2978.46)"
$ws.Range("E11").Value = "landscape_mesh"

# Row 12: landscape_np
$ws.Range("A12").Value = "Number of Patches"
$ws.Range("B12").Value = "340 (340–2983)"
$ws.Range("C12").Value = "809.76 (1–16279)"
$ws.Range("E12").Value = "landscape_np"

# Row 13: landscape_shei
$ws.Range("A13").Value = "Shannon's Evenness"
$ws.Range("B13").Value = "0.64 (0.64–0.87)"
$ws.Range("C13").Value = "0.62 (0–1)"
$ws.Range("E13").Value = "landscape_shei"

# Row 14: landscape_siei
$ws.Range("A14").Value = "Simpson's Evenness"
$ws.Range("B14").Value = "0.75 (0.75–0.9)"
$ws.Range("C14").Value = "0.7 (0–1)"
$ws.Range("E14").Value = "landscape_siei"

# Row 15: landscape_tca
$ws.Range("A15").Value = "Total Core Area"
$ws.Range("B15").Value = "0 (0–0)"
$ws.Range("C15").Value = "0 (0–0)"
$ws.Range("E15").Value = "landscape_tca"

# Row 16: lc_broadleaf
$ws.Range("A16").Value = "Broadleaf"
$ws.Range("B16").Value = "0.08 (0.08–0.52)"
$ws.Range("C16").Value = "0.15 (0–1)"
$ws.Range("E16").Value = "lc_broadleaf"

# Row 17: lc_coniferous
$ws.Range("A17").Value = "Coniferous"
$ws.Range("B17").Value = "0.34 (0.34–0.67)"
$ws.Range("C17").Value = "0.35 (0–1)"
$ws.Range("E17").Value = "lc_coniferous"

# Row 18: lc_herbs
$ws.Range("A18").Value = "Herbs"
$ws.Range("B18").Value = "0.01 (0.01–0.13)"
$ws.Range("C18").Value = "0.03 (0–0.43)"
$ws.Range("E18").Value = "lc_herbs"

# Row 19: lc_mixedwood
$ws.Range("A19").Value = "Mixedwood"
$ws.Range("B19").Value = "0.05 (0.05–0.21)"
$ws.Range("C19").Value = "0.07 (0–0.71)"
$ws.Range("E19").Value = "lc_mixedwood"

# Row 20: lc_shrubs
$ws.Range("A20").Value = "Shrubland"
$ws.Range("B20").Value = "0 (0–0.2)"
$ws.Range("C20").Value = "0.03 (0–0.96)"
$ws.Range("E20").Value = "lc_shrubs"

# Row 21: lc_wetland
$ws.Range("A21").Value = "Wetland"
$ws.Range("B21").Value = "0.04 (0.04–0.24)"
$ws.Range("C21").Value = "0.07 (0–0.76)"
$ws.Range("E21").Value = "lc_wetland"

# Row 22: lc_wetland_treed
$ws.Range("A22").Value = "Treed Wetland"
$ws.Range("B22").Value = "0.23 (0.23–0.6)"
$ws.Range("C22").Value = "0.25 (0–0.95)"
$ws.Range("E22").Value = "lc_wetland_treed"

# Row 23: nonanthro_cai_mn
$ws.Range("A23").Value = "Core Area Index (natural habitat)"
$ws.Range("B23").Value = "32.73 (32.73–84.77)"
$ws.Range("C23").Value = "36.34 (0–98.21)"
$ws.Range("E23").Value = "nonanthro_cai_mn"

# Row 24: nonanthro_ed
$ws.Range("A24").Value = "Edge Density"
$ws.Range("B24").Value = "74.6 (74.6–296.57)"
$ws.Range("C24").Value = "101.65 (0–742.03)"
$ws.Range("E24").Value = "nonanthro_ed"

# Row 25: nonanthro_tca
$ws.Range("A25").Value = "Total Core Area (natural)"
$ws.Range("B25").Value = "1024.54 (1024.54–5263.56)"
$ws.Range("C25").Value = "1664.25 (0–7585.52)"
$ws.Range("E25").Value = "nonanthro_tca"

# Row 26: osm_industrial
$ws.Range("A26").Value = "Industrial Facilities"
$ws.Range("B26").Value = "0 (0–0.09)"
$ws.Range("C26").Value = "0.02 (0–1)"
$ws.Range("E26").Value = "osm_industrial"

# Row 27: pct_lari_lar
$ws.Range("A27").Value = "Tamarack"
$ws.Range("B27").Value = "0 (0–0.09)"
$ws.Range("C27").Value = "0.02 (0–0.47)"
$ws.Range("E27").Value = "pct_lari_lar"

# Row 28: pct_pice_gla
$ws.Range("A28").Value = "White Spruce"
$ws.Range("B28").Value = "0 (0–0.06)"
$ws.Range("C28").Value = "0.01 (0–0.39)"
$ws.Range("E28").Value = "pct_pice_gla"

# Row 29: pct_pice_mar
$ws.Range("A29").Value = "Black Spruce"
$ws.Range("B29").Value = "0.56 (0.56–0.9)"
$ws.Range("C29").Value = "0.52 (0–1)"
$ws.Range("E29").Value = "pct_pice_mar"

# Row 30: pct_pinu_ban
$ws.Range("A30").Value = "Jack Pine"
$ws.Range("B30").Value = "0.01 (0.01–0.55)"
$ws.Range("C30").Value = "0.07 (0–1)"
$ws.Range("E30").Value = "pct_pinu_ban"

# Row 31: pct_popu_tre
$ws.Range("A31").Value = "Trembling Aspen"
$ws.Range("B31").Value = "0.32 (0.32–0.8)"
$ws.Range("C31").Value = "0.36 (0–1)"
$ws.Range("E31").Value = "pct_popu_tre"

# Row 32: pipe_trans
$ws.Range("A32").Value = "Pipelines & Transmission Lines"
$ws.Range("B32").Value = "0.01 (0.01–0.07)"
$ws.Range("C32").Value = "0.02 (0–0.46)"
$ws.Range("E32").Value = "pipe_trans"

# Row 33: roads
$ws.Range("A33").Value = "Roads"
$ws.Range("B33").Value = "0 (0–0.02)"
$ws.Range("C33").Value = "0.01 (0–0.12)"
$ws.Range("E33").Value = "roads"

# Row 34: seismic
$ws.Range("A34").Value = "Seismic Lines"
$ws.Range("B34").Value = "0.01 (0.01–0.04)"
$ws.Range("C34").Value = "0.01 (0–0.2)"
$ws.Range("E34").Value = "seismic"

# Row 35: wells_active
$ws.Range("A35").Value = "Active Well Sites"
$ws.Range("B35").Value = "0 (0–0.03)"
$ws.Range("C35").Value = "0.01 (0–0.32)"
$ws.Range("E35").Value = "wells_active"

# Row 36: wells_inactive
$ws.Range("A36").Value = "Inactive Well Sites"
$ws.Range("B36").Value = "0 (0–0.02)"
$ws.Range("C36").Value = "0.01 (0–0.59)"
$ws.Range("E36").Value = "wells_inactive"

# Row 37: wells_total
$ws.Range("A37").Value = "Well Sites"
$ws.Range("B37").Value = "0 (0–0.05)"
$ws.Range("C37").Value = "0.01 (0–0.59)"
$ws.Range("E37").Value = "wells_total"

